# Add two new weekly price rows for "Membrillo" (Vega Central Mapocho de Santiago)
# by inserting them at row 49 - this shifts the existing rows 49:75 down to 51:77.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (row 49)
$ws.Rows.Item(49).Insert()
$ws.Rows.Item(49).Insert()

# ---- New row 49 ----
$ws.Cells.Item(49, 1).Value = 9
$ws.Cells.Item(49, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(49, 3).Value = "Metropolitana"
$ws.Cells.Item(49, 4).Value = 45062
$ws.Cells.Item(49, 5).Value = 13
$ws.Cells.Item(49, 6).Value = "Fruta"
$ws.Cells.Item(49, 7).Value = 100104
$ws.Cells.Item(49, 8).Value = "Frutos de pepita"
$ws.Cells.Item(49, 9).Value = 100104003
$ws.Cells.Item(49, 10).Value = "Membrillo"
$ws.Cells.Item(49, 11).Value = "Champion"
$ws.Cells.Item(49, 12).Value = "Especial"
$ws.Cells.Item(49, 13).Value = 280
$ws.Cells.Item(49, 14).Value = 10000
$ws.Cells.Item(49, 15).Value = 10500
$ws.Cells.Item(49, 16).Value = 10268
$ws.Cells.Item(49, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(49, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(49, 19).Value = 685
$ws.Cells.Item(49, 20).Value = 15

# ---- New row 50 ----
$ws.Cells.Item(50, 1).Value = 9
$ws.Cells.Item(50, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(50, 3).Value = "Metropolitana"
$ws.Cells.Item(50, 4).Value = 45062
$ws.Cells.Item(50, 5).Value = 13
$ws.Cells.Item(50, 6).Value = "Fruta"
$ws.Cells.Item(50, 7).Value = 100104
$ws.Cells.Item(50, 8).Value = "Frutos de pepita"
$ws.Cells.Item(50, 9).Value = 100104003
$ws.Cells.Item(50, 10).Value = "Membrillo"
$ws.Cells.Item(50, 11).Value = "Champion"
$ws.Cells.Item(50, 12).Value = "Primera"
$ws.Cells.Item(50, 13).Value = 350
$ws.Cells.Item(50, 14).Value = 7000
$ws.Cells.Item(50, 15).Value = 7500
$ws.Cells.Item(50, 16).Value = 7286
$ws.Cells.Item(50, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(50, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(50, 19).Value = 486
$ws.Cells.Item(50, 20).Value = 15
